# Applies the "Correcciones documentación Tablas y Relaciones" fixes:
#  - Merge the split "APROBADO  (BIN" + ")" runs into a single run "APROBADO  (BIN)"
#    and drop the stray trailing <a:endParaRPr/> on that paragraph.
#  - Drop stray trailing <a:endParaRPr/> elements on several other table-cell
#    paragraphs (text content itself is unchanged for these).
#
# All of the affected paragraphs live in "big relational table" shapes on
# slide 1, in the first column (column 1) of the third table row.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Fix-CellParagraph {
    param(
        [int]$ShapeIndex,
        [int]$Row,
        [int]$Col,
        [int]$ParaIndex
    )

    $sh = $s.Shapes.Item($ShapeIndex)
    $tbl = $sh.Table
    $cell = $tbl.Cell($Row, $Col)
    $tr = $cell.Shape.TextFrame.TextRange
    $paragraphs = $tr.Paragraphs()
    $para = $paragraphs.Item($ParaIndex)
    # Re-assigning the paragraph's own text normalizes its run/endParaRPr
    # structure: any split runs collapse into one run, and a stray trailing
    # endParaRPr (left over from earlier manual edits in PowerPoint) is
    # dropped.
    $para.Text = $para.Text
}

# Shape 2 ("2 Tabla"): merge "APROBADO  (BIN" + ")" into one run and drop
# the trailing endParaRPr.
Fix-CellParagraph -ShapeIndex 2 -Row 3 -Col 1 -ParaIndex 4
# Shape 2 ("2 Tabla"): ID_JP
Fix-CellParagraph -ShapeIndex 2 -Row 3 -Col 1 -ParaIndex 7

# Shape 3 ("3 Tabla"): INICIO_TRAMITACION, TIPO_CONTRATACION, TRAM_ANTICIPADA
Fix-CellParagraph -ShapeIndex 3 -Row 3 -Col 1 -ParaIndex 9
Fix-CellParagraph -ShapeIndex 3 -Row 3 -Col 1 -ParaIndex 10
Fix-CellParagraph -ShapeIndex 3 -Row 3 -Col 1 -ParaIndex 19

# Shape 4 ("4 Tabla"): ID_JP
Fix-CellParagraph -ShapeIndex 4 -Row 3 -Col 1 -ParaIndex 2

# Shape 6 ("8 Tabla"): PERIODO
Fix-CellParagraph -ShapeIndex 6 -Row 3 -Col 1 -ParaIndex 7

# Shape 8 ("12 Tabla"): TIPO_EVENTO_CONTRATACION
Fix-CellParagraph -ShapeIndex 8 -Row 3 -Col 1 -ParaIndex 2

# Shape 43 ("50 Tabla"): TIPO_ATRIBUTO
Fix-CellParagraph -ShapeIndex 43 -Row 3 -Col 1 -ParaIndex 1

# Shape 49 ("72 Tabla"): NUM_FACTURA
Fix-CellParagraph -ShapeIndex 49 -Row 3 -Col 1 -ParaIndex 1

# Shape 52 ("75 Tabla"): SP, TIPO_EVENTO_PLANIFICACION
Fix-CellParagraph -ShapeIndex 52 -Row 3 -Col 1 -ParaIndex 1
Fix-CellParagraph -ShapeIndex 52 -Row 3 -Col 1 -ParaIndex 2
